# Horarios actualizados Linea 141 - 134
# Applies the scraped update (new "Ultima actualizacion" timestamp 07:12:46,
# and new rows of data) to the three worksheets of the workbook.

$wb = $excel.ActiveWorkbook

$updateTime = "Última actualización: 07:12:46"

# ---------------------------------------------------------------------------
# Sheet 1: LP1912  -> grows from 39 to 42 data rows (rows 6-44 -> 6-47)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = $updateTime
$ws1.Range("A3").Value = "Total filas: 42"

$ws1.Cells.Item(45, 1).Value = "07:12:46"
$ws1.Cells.Item(45, 2).Value = "08:51"
$ws1.Cells.Item(45, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(45, 4).Value = 99
$ws1.Cells.Item(45, 5).Value = "LP1912"

$ws1.Cells.Item(46, 1).Value = "07:12:46"
$ws1.Cells.Item(46, 2).Value = "08:53"
$ws1.Cells.Item(46, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(46, 4).Value = 101
$ws1.Cells.Item(46, 5).Value = "LP1912"

$ws1.Cells.Item(47, 1).Value = "07:12:46"
$ws1.Cells.Item(47, 2).Value = "08:57"
$ws1.Cells.Item(47, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(47, 4).Value = 105
$ws1.Cells.Item(47, 5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215 -> grows from 10 to 12 data rows (rows 6-15 -> 6-17)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = $updateTime
$ws2.Range("A3").Value = "Total filas: 12"

$ws2.Cells.Item(16, 1).Value = "07:12:46"
$ws2.Cells.Item(16, 2).Value = "08:53"
$ws2.Cells.Item(16, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(16, 4).Value = 101
$ws2.Cells.Item(16, 5).Value = "LP1912"

$ws2.Cells.Item(17, 1).Value = "07:12:46"
$ws2.Cells.Item(17, 2).Value = "08:57"
$ws2.Cells.Item(17, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(17, 4).Value = 105
$ws2.Cells.Item(17, 5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173 -> grows from 5 to 6 data rows (rows 6-10 -> 6-11)
# A new row is inserted at row 10, pushing the previous row 10 down to row 11
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = $updateTime
$ws3.Range("A3").Value = "Total filas: 6"

$ws3.Rows.Item(10).Insert()

$ws3.Cells.Item(10, 1).Value = "07:12:46"
$ws3.Cells.Item(10, 2).Value = "08:50"
$ws3.Cells.Item(10, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(10, 4).Value = 98
$ws3.Cells.Item(10, 5).Value = "L6203"
